$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 1: Exercise_1_Data  -- fill in week 1 field data for 3 trees
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Exercise_1_Data")

# Tree 1 - eastern white pine (Pinus strobus)
$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = "eastern white pine (Pinus strobus)"
$ws1.Range("C2").Value = "DBH"
$ws1.Range("D2").Value = "Diameter tape"
$ws1.Range("E2").Value = 37.2

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "eastern white pine (Pinus strobus)"
$ws1.Range("C3").Value = "DBH"
$ws1.Range("D3").Value = "Calipers"
$ws1.Range("E3").Value = 35.6

$ws1.Range("A4").Value = 1
$ws1.Range("B4").Value = "eastern white pine (Pinus strobus)"
$ws1.Range("C4").Value = "Height"
$ws1.Range("D4").Value = "Clinometer"
$ws1.Range("E4").Value = 25.65

$ws1.Range("A5").Value = 1
$ws1.Range("B5").Value = "eastern white pine (Pinus strobus)"
$ws1.Range("C5").Value = "Height"
$ws1.Range("D5").Value = "Staff hypsometer"
$ws1.Range("E5").Value = 26

$ws1.Range("A6").Value = 1
$ws1.Range("B6").Value = "eastern white pine (Pinus strobus)"
$ws1.Range("C6").Value = "Height"
$ws1.Range("D6").Value = "Staff hypsometer"

# Tree 2 - black birch (Betula lenta)
$ws1.Range("A7").Value = 2
$ws1.Range("B7").Value = "black birch (Betula lenta)"
$ws1.Range("C7").Value = "Height"
$ws1.Range("D7").Value = "Clinometer"
$ws1.Range("E7").Value = 24.92

$ws1.Range("A8").Value = 2
$ws1.Range("B8").Value = "black birch (Betula lenta)"
$ws1.Range("C8").Value = "DBH"
$ws1.Range("D8").Value = "Diameter tape"
$ws1.Range("E8").Value = 33

$ws1.Range("A9").Value = 2
$ws1.Range("B9").Value = "black birch (Betula lenta)"
$ws1.Range("C9").Value = "DBH"
$ws1.Range("D9").Value = "Calipers"
$ws1.Range("E9").Value = 31

$ws1.Range("A10").Value = 2
$ws1.Range("B10").Value = "black birch (Betula lenta)"
$ws1.Range("C10").Value = "Height"
$ws1.Range("D10").Value = "Biltmore stick"

$ws1.Range("A11").Value = 2
$ws1.Range("B11").Value = "black birch (Betula lenta)"
$ws1.Range("C11").Value = "Height"
$ws1.Range("D11").Value = "Staff hypsometer"

# Tree 3 - red oak (Quercus rubra)
$ws1.Range("A12").Value = 3
$ws1.Range("B12").Value = "red oak (Quercus rubra)"
$ws1.Range("C12").Value = "DBH"
$ws1.Range("D12").Value = "Biltmore stick"
$ws1.Range("E12").Value = 40

$ws1.Range("A13").Value = 3
$ws1.Range("B13").Value = "red oak (Quercus rubra)"
$ws1.Range("C13").Value = "DBH"
$ws1.Range("D13").Value = "Diameter tape"
$ws1.Range("E13").Value = 41

$ws1.Range("A14").Value = 3
$ws1.Range("B14").Value = "red oak (Quercus rubra)"
$ws1.Range("C14").Value = "DBH"
$ws1.Range("D14").Value = "Calipers"
$ws1.Range("E14").Value = 38.1

$ws1.Range("A15").Value = 3
$ws1.Range("B15").Value = "red oak (Quercus rubra)"
$ws1.Range("C15").Value = "Height"
$ws1.Range("D15").Value = "Staff hypsometer"
$ws1.Range("E15").Value = 19

$ws1.Range("A16").Value = 3
$ws1.Range("B16").Value = "red oak (Quercus rubra)"
$ws1.Range("C16").Value = "Height"
$ws1.Range("D16").Value = "Clinometer"
$ws1.Range("E16").Value = 16.4

$ws1.Range("E16").Select()
$ws1.Application.ActiveWindow.ScrollColumn = 2

# ------------------------------------------------------------------
# Sheet 2: Clinometer_Calculations -- clinometer readings for trees
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Clinometer_Calculations")

$ws2.Range("A2").Value = 2
$ws2.Range("C2").Value = 17.8
$ws2.Range("D2").Value = 150
$ws2.Range("E2").Value = 10

$ws2.Range("C3").Value = 19
$ws2.Range("D3").Value = 140
$ws2.Range("E3").Value = 5

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = "red oak (Quercus rubra)"
$ws2.Range("C4").Value = 28
$ws2.Range("D4").Value = 72
$ws2.Range("E4").Value = 14

$ws2.Range("D37").Select()
$ws2.Application.ActiveWindow.ScrollColumn = 3

# ------------------------------------------------------------------
# Sheet 3: Field_Source -- add blank marker value for instrument list
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Field_Source")
$ws3.Range("D7").Value = " "

$ws3.Range("D7").Select()
$ws3.Application.ActiveWindow.ScrollRow = 6

$ws1.Activate()
